$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.545.23'
$ws.Range("E2").Value = '  +0.27%  '

$ws.Range("D3").Value = '1.577.66'
$ws.Range("E3").Value = '  +0.38%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '289.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.69%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3700'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.69'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.55%  '

$ws.Range("E9").Value = '  -0.78%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.143'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07487'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.01'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.49%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.005'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.954'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.22%  '

$ws.Range("D16").Value = '1.581.50'
$ws.Range("E16").Value = '  +0.64%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001122'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.36%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '88.66'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.00%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06769'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.14%  '

$ws.Range("E20").Value = '  +0.05%  '

$ws.Range("E21").Value = '  +1.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.58'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.18'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.53%  '

$ws.Range("D24").Value = '22.531.80'
$ws.Range("E24").Value = '  +0.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.406'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.603'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.84'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.52%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.53%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.021'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.54%  '

$ws.Range("D31").Value = '1.756.90'
$ws.Range("E31").Value = '  +0.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.064'
$ws.Range("D32").Style = "Normal"

$ws.Range("E33").Value = '  -0.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.009'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.691'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.60%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08338'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.14%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02461'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.62%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2303'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.441'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06394'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.96%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.299'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6357'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.48%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.77%  '

$ws.Range("E44").Value = '  +0.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.98'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6204'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.773'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.066'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.77'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.218'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07273'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.61%  '
